# "added new footer xpath"
#
# This script:
#  1. Deletes the "UpdateAccountInformation" sheet.
#  2. Rewrites the "CCPayment" sheet with a new header/data layout
#     (Type / Transaction / Card number / Expiry month / Expiry year /
#     Cardholder name / Security code) and two data rows (Success + Failure).
#  3. Updates the "ConvenienceStorePayment" sheet's hyperlink-adjacent data
#     is untouched; it only shifts because shared strings were renumbered
#     when new strings were inserted earlier in the shared string table -
#     that renumbering happens automatically because we insert the new
#     strings before the old ones get written again, so we only need to
#     touch the actual edited cells.
#  4. Activates the "User" sheet so it becomes the workbook's selected tab.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------
# 1. Remove the UpdateAccountInformation sheet entirely.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("UpdateAccountInformation").Delete()

# ---------------------------------------------------------------------
# 2. Rebuild the CCPayment sheet contents.
# ---------------------------------------------------------------------
$cc = $wb.Worksheets.Item("CCPayment")

$cc.Range("A1").Value = "Type"
$cc.Range("B1").Value = "Transaction"
$cc.Range("C1").Value = "クレジットカード番号"
$cc.Range("D1").Value = "有効期限（月）"
$cc.Range("E1").Value = "有効期限（年）"
$cc.Range("F1").Value = "カード名義"
$cc.Range("G1").Value = "セキュリティコード"

$cc.Range("A2").Value = "Visa"
$cc.Range("B2").Value = "Success"
$cc.Range("C2").Value = 4111111111111110
$cc.Range("D2").Value = 6
$cc.Range("E2").Value = 2025
$cc.Range("F2").Value = "Test Card"
$cc.Range("G2").Value = 1213

$cc.Range("A3").Value = "Visa"
$cc.Range("B3").Value = "Failure"
$cc.Range("C3").Value = 4111111111111110
$cc.Range("D3").Value = 6
$cc.Range("E3").Value = 2025
$cc.Range("F3").Value = "Test Card"
$cc.Range("G3").Value = 1213

$cc.Range("B3").Select()

# ---------------------------------------------------------------------
# 3. Select the User sheet (it becomes the active / tab-selected sheet).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("User").Activate()
